$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '69.566.44'
$ws.Range('E2').Value = '  +0.75%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.774.94'
$ws.Range('E3').Value = '  -0.14%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.03%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '660.05'
$ws.Range('E5').Value = '  +4.98%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '165.74'
$ws.Range('E6').Value = '  +1.00%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '3.773.28'
$ws.Range('E7').Value = '  -0.08%  '

$ws.Range('E8').Value = '  +0.03%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.526'
$ws.Range('E9').Value = '  +1.27%  '

$ws.Range('E10').Value = '  -0.68%  '

$ws.Range('E11').Value = '  +1.35%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '6.98'
$ws.Range('E12').Value = '  +5.14%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000240'
$ws.Range('E13').Value = '  -2.79%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '4.414.20'
$ws.Range('E15').Value = '  +0.08%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '3.767.90'
$ws.Range('E16').Value = '  +1.56%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '69.520.76'
$ws.Range('E17').Value = '  +0.85%  '

$ws.Range('E18').Value = '  -1.40%  '

$ws.Range('E19').Value = '  +0.74%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '7.11'
$ws.Range('E20').Value = '  +0.31%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '469.42'
$ws.Range('E21').Value = '  +0.38%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '9.88'
$ws.Range('E22').Value = '  +2.46%  '

$ws.Range('E23').Value = '  +1.06%  '

$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '82.30'
$ws.Range('E24').Value = '  -1.03%  '

$ws.Range('B25').Value = 'PEPE'
$ws.Range('C25').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.0000144'
$ws.Range('E25').Value = '  -3.82%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '12.25'
$ws.Range('E26').Value = '  +1.93%  '

$ws.Range('E27').Value = '  +2.84%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.11'
$ws.Range('E28').Value = '  -1.50%  '

$ws.Range('E29').Value = '  +0.16%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '3.926.20'
$ws.Range('E30').Value = '  +0.05%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '2.76'
$ws.Range('E31').Value = '  +3.23%  '

$ws.Range('E32').Value = '  +3.57%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '7.28'
$ws.Range('E33').Value = '  +0.22%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '28.90'
$ws.Range('E34').Value = '  -0.14%  '

$ws.Range('E35').Value = '  +18.56%  '

$ws.Range('E36').Value = '  -0.11%  '

$ws.Range('B37').Value = 'Aptos'
$ws.Range('C37').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '9.00'

$ws.Range('B38').Value = 'RenzoRestakedETH'
$ws.Range('C38').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.732.46'
$ws.Range('E38').Value = '  +0.35%  '

$ws.Range('E40').Value = '  +0.36%  '

$ws.Range('E41').Value = '  -1.44%  '

$ws.Range('E42').Value = '  +0.27%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.960'
$ws.Range('E43').Value = '  -1.04%  '

$ws.Range('E44').Value = '  -0.05%  '

$ws.Range('B45').Value = 'Stacks'
$ws.Range('C45').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.03'
$ws.Range('E45').Value = '  +6.14%  '

$ws.Range('B46').Value = 'Arweave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '45.28'
$ws.Range('E46').Value = '  +6.22%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '157.66'
$ws.Range('E47').Value = '  +2.76%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '48.04'
$ws.Range('E48').Value = '  +2.50%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.298'
$ws.Range('E49').Value = '  +0.32%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.39'
$ws.Range('E50').Value = '  +1.02%  '

$ws.Range('E51').Value = '  +0.58%  '
